# Append a new measurement row (row 8) to the temperature history sheet,
# matching the format of the existing data rows (6 and 7): a timestamp
# string, a numeric temperature, and a humidity percentage stored as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A8: timestamp, stored as plain text (same as A6/A7).
$ws.Range("A8").Value = "28/03/2025 17:03:18"

# B8: temperature, stored as a number (same as B6/B7).
$ws.Range("B8").Value = 20.2

# C8: humidity, stored as literal text "94%" (same as C6/C7), not as a
# numeric percentage. Excel normally auto-converts a "94%"-looking value
# into a percentage number, so force text interpretation via the "@"
# number format while assigning it, then drop back to the default/Normal
# style so the cell ends up unstyled, just like the other data cells.
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "94%"
$ws.Range("C8").Style = "Normal"
